$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text format so numeric-looking strings (thousand-dot
# separators, trailing zeros, scientific-looking decimals) are preserved exactly
# as text instead of being reinterpreted by Excel as numbers.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

# Volume(1h) column (E) values already contain non-numeric characters (%, spaces)
# so they round-trip as text without any special handling.
function Set-PercentValue($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

Set-TextValue 'D2' '27.781.66'
Set-PercentValue 'E2' '  +1.69%  '
Set-TextValue 'D3' '1.885.68'
Set-PercentValue 'E3' '  +1.74%  '
Set-PercentValue 'E4' '  +0.37%  '
Set-TextValue 'D5' '333.79'
Set-PercentValue 'E5' '  +1.61%  '
Set-PercentValue 'E6' '  +0.40%  '
Set-TextValue 'D7' '0.4719'
Set-PercentValue 'E7' '  +2.42%  '
Set-PercentValue 'E8' '  -0.45%  '
Set-TextValue 'D9' '47.64'
Set-PercentValue 'E9' '  +1.99%  '
Set-TextValue 'D10' '0.08068'
Set-PercentValue 'E10' '  +1.67%  '
Set-TextValue 'D11' '1.027'
Set-PercentValue 'E11' '  +1.62%  '
Set-TextValue 'D12' '22.14'
Set-PercentValue 'E12' '  +3.18%  '
Set-TextValue 'D13' '1.884.19'
Set-PercentValue 'E13' '  +2.29%  '
Set-TextValue 'D14' '5.975'
Set-PercentValue 'E14' '  +0.92%  '
Set-TextValue 'D15' '7.125'
Set-PercentValue 'E15' '  -0.08%  '
Set-TextValue 'D16' '1.010'
Set-PercentValue 'E16' '  +0.66%  '
Set-TextValue 'D17' '0.06755'
Set-PercentValue 'E17' '  +2.57%  '
Set-TextValue 'D18' '87.24'
Set-PercentValue 'E18' '  +1.30%  '
Set-TextValue 'D19' '0.00001047'
Set-PercentValue 'E19' '  +1.66%  '
Set-TextValue 'D20' '17.35'
Set-PercentValue 'E20' '  +0.72%  '
Set-PercentValue 'E21' '  +0.45%  '
Set-TextValue 'D22' '27.805.08'
Set-PercentValue 'E22' '  +1.75%  '
Set-TextValue 'D23' '5.521'
Set-PercentValue 'E23' '  +0.95%  '
Set-TextValue 'D24' '11.00'
Set-PercentValue 'E24' '  +1.12%  '
Set-TextValue 'D25' '2.332'
Set-PercentValue 'E25' '  +1.54%  '
Set-TextValue 'D26' '2.108.32'
Set-PercentValue 'E26' '  +2.12%  '
Set-TextValue 'D27' '158.96'
Set-PercentValue 'E27' '  +3.64%  '
Set-PercentValue 'E28' '  -0.13%  '
Set-TextValue 'D29' '2.104'
Set-PercentValue 'E29' '  +2.17%  '
Set-TextValue 'D30' '5.568'
Set-PercentValue 'E30' '  +2.20%  '
Set-TextValue 'D31' '122.04'
Set-PercentValue 'E31' '  +0.44%  '
Set-TextValue 'D32' '0.9809'
Set-PercentValue 'E32' '  +3.54%  '
Set-PercentValue 'E33' '  +0.78%  '
Set-PercentValue 'E34' '  +0.82%  '
Set-TextValue 'D35' '3.616'
Set-PercentValue 'E35' '  +0.81%  '
Set-TextValue 'D36' '5.352'
Set-PercentValue 'E36' '  +1.77%  '
Set-TextValue 'D37' '0.06154'
Set-PercentValue 'E37' '  +2.07%  '
Set-TextValue 'D38' '0.02263'
Set-PercentValue 'E38' '  +1.64%  '
Set-TextValue 'D39' '1.218'
Set-PercentValue 'E39' '  +0.76%  '
Set-TextValue 'D40' '8.071'
Set-PercentValue 'E40' '  +0.57%  '
Set-TextValue 'D41' '0.6001'
Set-PercentValue 'E41' '  +1.43%  '
Set-TextValue 'D42' '0.1894'
Set-TextValue 'D43' '10.32'
Set-PercentValue 'E43' '  +1.63%  '
Set-PercentValue 'E44' '  -1.91%  '
Set-TextValue 'D45' '0.5708'
Set-PercentValue 'E45' '  +1.69%  '
Set-TextValue 'D46' '12.24'
Set-PercentValue 'E46' '  +1.50%  '
Set-TextValue 'D47' '1.944'
Set-PercentValue 'E47' '  +1.55%  '
Set-TextValue 'D48' '3.394'
Set-PercentValue 'E48' '  -0.06%  '
Set-TextValue 'D49' '0.06912'
Set-PercentValue 'E49' '  +2.22%  '
Set-TextValue 'D50' '114.13'
Set-PercentValue 'E50' '  +4.76%  '
Set-TextValue 'D51' '0.00000000306'
Set-PercentValue 'E51' '  -0.86%  '
